$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 53 (Leve Item ID = 5479)
$ws.Range("H53").Value = 220.8
$ws.Range("J53").Value = 209.54546
$ws.Range("L53").Value = 209.54546
$ws.Range("N53").Value = -1483.54546
# row 112 (Leve Item ID = 27960)
$ws.Range("H112").Value = 1684.725
$ws.Range("I112").Value = 775
$ws.Range("J112").Value = 1785.8055
$ws.Range("K112").Value = 2325
$ws.Range("L112").Value = 5357.416499999999
$ws.Range("M112").Value = -1217
$ws.Range("N112").Value = -7573.416499999999
# row 125 (Leve Item ID = 36228)
$ws.Range("H125").Value = 701.7619
$ws.Range("I125").Value = 490.41666
$ws.Range("J125").Value = 983.55554
$ws.Range("K125").Value = 4413.74994
$ws.Range("L125").Value = 8851.99986
$ws.Range("M125").Value = -1953.74994
$ws.Range("N125").Value = -13771.99986
# row 129 (Leve Item ID = 36115)
$ws.Range("H129").Value = 1608.5238
$ws.Range("I129").Value = 476.46155
$ws.Range("J129").Value = 3448.125
$ws.Range("K129").Value = 1429.38465
$ws.Range("L129").Value = 10344.375
$ws.Range("M129").Value = 3570.61535
$ws.Range("N129").Value = -20344.375

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 63 (Leve Item ID = 12528)
$ws.Range("H63").Value = 5693.4443
$ws.Range("I63").Value = 1748.7142
$ws.Range("J63").Value = 19500
$ws.Range("K63").Value = 1748.7142
$ws.Range("L63").Value = 19500
$ws.Range("M63").Value = -1062.7142
$ws.Range("N63").Value = -20872
# row 66 (Leve Item ID = 12528)
$ws.Range("H66").Value = 5693.4443
$ws.Range("I66").Value = 1748.7142
$ws.Range("J66").Value = 19500
$ws.Range("K66").Value = 8743.571
$ws.Range("L66").Value = 97500
$ws.Range("M66").Value = -5311.571
$ws.Range("N66").Value = -104364
# row 74 (Leve Item ID = 44000)
$ws.Range("H74").Value = 1760.4546
$ws.Range("I74").Value = 1800.258
$ws.Range("J74").Value = 1665.5385
$ws.Range("K74").Value = 1800.258
$ws.Range("L74").Value = 1665.5385
$ws.Range("M74").Value = -926.258
$ws.Range("N74").Value = -3413.5385
# row 77 (Leve Item ID = 44000)
$ws.Range("H77").Value = 1760.4546
$ws.Range("I77").Value = 1800.258
$ws.Range("J77").Value = 1665.5385
$ws.Range("K77").Value = 9001.290000000001
$ws.Range("L77").Value = 8327.692500000001
$ws.Range("M77").Value = -4633.290000000001
$ws.Range("N77").Value = -17063.6925
# row 103 (Leve Item ID = 18533)
$ws.Range("H103").Value = 25000
$ws.Range("J103").Value = 25000
$ws.Range("L103").Value = 25000
$ws.Range("N103").Value = -27344
# row 107 (Leve Item ID = 25645)
$ws.Range("H107").Value = 90000
$ws.Range("J107").Value = 90000
$ws.Range("L107").Value = 90000
$ws.Range("N107").Value = -97680

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 82 (Leve Item ID = 11877)
$ws.Range("H82").Value = 10841.546
$ws.Range("I82").Value = 5695.222
$ws.Range("J82").Value = 34000
$ws.Range("K82").Value = 5695.222
$ws.Range("L82").Value = 34000
$ws.Range("M82").Value = -5312.222
$ws.Range("N82").Value = -34766
# row 85 (Leve Item ID = 11877)
$ws.Range("H85").Value = 10841.546
$ws.Range("I85").Value = 5695.222
$ws.Range("J85").Value = 34000
$ws.Range("K85").Value = 5695.222
$ws.Range("L85").Value = 34000
$ws.Range("M85").Value = -4369.222
$ws.Range("N85").Value = -36652

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID = 44023)
$ws.Range("H31").Value = 2223.1614
$ws.Range("I31").Value = 1317.579
$ws.Range("J31").Value = 3657
$ws.Range("K31").Value = 1317.579
$ws.Range("L31").Value = 3657
$ws.Range("M31").Value = -1022.579
$ws.Range("N31").Value = -4247
# row 34 (Leve Item ID = 44023)
$ws.Range("H34").Value = 2223.1614
$ws.Range("I34").Value = 1317.579
$ws.Range("J34").Value = 3657
$ws.Range("K34").Value = 1317.579
$ws.Range("L34").Value = 3657
$ws.Range("M34").Value = -1115.579
$ws.Range("N34").Value = -4061
# row 58 (Leve Item ID = 44021)
$ws.Range("H58").Value = 1424.8269
$ws.Range("I58").Value = 805.32556
$ws.Range("J58").Value = 4384.6665
$ws.Range("K58").Value = 805.32556
$ws.Range("L58").Value = 4384.6665
$ws.Range("M58").Value = -602.32556
$ws.Range("N58").Value = -4790.6665
# row 134 (Leve Item ID = 44020)
$ws.Range("H134").Value = 1555.85
$ws.Range("I134").Value = 1000.3774
$ws.Range("J134").Value = 5761.5713
$ws.Range("K134").Value = 3001.1322
$ws.Range("L134").Value = 17284.7139
$ws.Range("M134").Value = -466.1322
$ws.Range("N134").Value = -22354.7139
# row 136 (Leve Item ID = 44021)
$ws.Range("H136").Value = 1424.8269
$ws.Range("I136").Value = 805.32556
$ws.Range("J136").Value = 4384.6665
$ws.Range("K136").Value = 2415.97668
$ws.Range("L136").Value = 13153.9995
$ws.Range("M136").Value = 134.0233200000002
$ws.Range("N136").Value = -18253.9995

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 46 (Leve Item ID = 4701)
$ws.Range("H46").Value = 542
$ws.Range("I46").Value = 562.5
$ws.Range("J46").Value = 501
$ws.Range("K46").Value = 1687.5
$ws.Range("L46").Value = 1503
$ws.Range("M46").Value = -1596.5
$ws.Range("N46").Value = -1685
# row 92 (Leve Item ID = 19841)
$ws.Range("H92").Value = 725
$ws.Range("J92").Value = 725
$ws.Range("L92").Value = 2175
$ws.Range("N92").Value = -4671
# row 109 (Leve Item ID = 27854)
$ws.Range("H109").Value = 3344.0232
$ws.Range("I109").Value = 674.9231
$ws.Range("J109").Value = 4500.6333
$ws.Range("K109").Value = 2024.7693
$ws.Range("L109").Value = 13501.8999
$ws.Range("M109").Value = -984.7692999999999
$ws.Range("N109").Value = -15581.8999
# row 113 (Leve Item ID = 27843)
$ws.Range("H113").Value = 4926566.5
$ws.Range("I113").Value = 5747543.5
$ws.Range("J113").Value = 704
$ws.Range("K113").Value = 17242630.5
$ws.Range("L113").Value = 2112
$ws.Range("M113").Value = -17240460.5
$ws.Range("N113").Value = -6452
# row 122 (Leve Item ID = 36078)
$ws.Range("H122").Value = 8565.629999999999
$ws.Range("J122").Value = 845.875
$ws.Range("L122").Value = 7612.875
$ws.Range("N122").Value = -12512.875
# row 124 (Leve Item ID = 36040)
$ws.Range("H124").Value = 2362.2222
$ws.Range("I124").Value = 1065
$ws.Range("J124").Value = 3400
$ws.Range("K124").Value = 3195
$ws.Range("L124").Value = 10200
$ws.Range("M124").Value = 1715
$ws.Range("N124").Value = -20020

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 10 (Leve Item ID = 4306)
$ws.Range("H10").Value = 250512500
$ws.Range("I10").Value = 250512500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 250512500
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -250512331
# row 28 (Leve Item ID = 2063)
$ws.Range("H28").Value = 3015
$ws.Range("J28").Value = 3015
$ws.Range("L28").Value = 3015
$ws.Range("N28").Value = -3399
# row 111 (Leve Item ID = 25853)
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -26134
# row 122 (Leve Item ID = 36182)
$ws.Range("H122").Value = 11111955
$ws.Range("I122").Value = 20000980
$ws.Range("J122").Value = 674.5
$ws.Range("K122").Value = 60002940
$ws.Range("L122").Value = 2023.5
$ws.Range("M122").Value = -60000490
$ws.Range("N122").Value = -6923.5
# row 126 (Leve Item ID = 36184)
$ws.Range("H126").Value = 4547562
$ws.Range("I126").Value = 7693764.5
$ws.Range("J126").Value = 3046.889
$ws.Range("K126").Value = 23081293.5
$ws.Range("L126").Value = 9140.667000000001
$ws.Range("M126").Value = -23078823.5
$ws.Range("N126").Value = -14080.667

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 46 (Leve Item ID = 5282)
$ws.Range("H46").Value = 999.6667
$ws.Range("I46").Value = 999.6667
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 999.6667
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -811.6667
# row 61 (Leve Item ID = 27740)
$ws.Range("H61").Value = 978.9545000000001
$ws.Range("I61").Value = 933.06665
$ws.Range("J61").Value = 1077.2858
$ws.Range("K61").Value = 933.06665
$ws.Range("L61").Value = 1077.2858
$ws.Range("M61").Value = -731.06665
$ws.Range("N61").Value = -1481.2858
# row 101 (Leve Item ID = 18549)
$ws.Range("H101").Value = 6240.5
$ws.Range("J101").Value = 6240.5
$ws.Range("L101").Value = 6240.5
$ws.Range("N101").Value = -12730.5
# row 113 (Leve Item ID = 27740)
$ws.Range("H113").Value = 978.9545000000001
$ws.Range("I113").Value = 933.06665
$ws.Range("J113").Value = 1077.2858
$ws.Range("K113").Value = 933.06665
$ws.Range("L113").Value = 1077.2858
$ws.Range("M113").Value = 1236.93335
$ws.Range("N113").Value = -5417.2858

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 64 (Leve Item ID = 11036)
$ws.Range("H64").Value = 12000
$ws.Range("J64").Value = 12000
$ws.Range("L64").Value = 12000
$ws.Range("N64").Value = -12496
# row 67 (Leve Item ID = 11036)
$ws.Range("H67").Value = 12000
$ws.Range("J67").Value = 12000
$ws.Range("L67").Value = 12000
$ws.Range("N67").Value = -13716
# row 98 (Leve Item ID = 18374)
$ws.Range("H98").Value = 18590
$ws.Range("J98").Value = 18590
$ws.Range("L98").Value = 18590
$ws.Range("N98").Value = -24580
# row 101 (Leve Item ID = 18538)
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490
# row 107 (Leve Item ID = 27746)
$ws.Range("H107").Value = 273.375
$ws.Range("I107").Value = 255.625
$ws.Range("J107").Value = 282.25
$ws.Range("K107").Value = 766.875
$ws.Range("L107").Value = 846.75
$ws.Range("M107").Value = 1153.125
$ws.Range("N107").Value = -4686.75
# row 113 (Leve Item ID = 27752)
$ws.Range("H113").Value = 550.86365
$ws.Range("I113").Value = 543.1053000000001
$ws.Range("K113").Value = 1629.3159
$ws.Range("M113").Value = 540.6840999999999
# row 132 (Leve Item ID = 44029)
$ws.Range("H132").Value = 1948.1136
$ws.Range("I132").Value = 1132.1305
$ws.Range("J132").Value = 2841.8096
$ws.Range("K132").Value = 3396.3915
$ws.Range("L132").Value = 8525.4288
$ws.Range("M132").Value = -866.3914999999997
$ws.Range("N132").Value = -13585.4288

